$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, pushing existing rows 42:91 down to 43:92
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly record
$ws.Cells.Item(42, 1).Value = 4
$ws.Cells.Item(42, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(42, 3).Value = "Los Lagos"
$ws.Cells.Item(42, 4).Value = 44494
$ws.Cells.Item(42, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(42, 5).Value = 10
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100108
$ws.Cells.Item(42, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(42, 9).Value = 100108002
$ws.Cells.Item(42, 10).Value = "Mango"
$ws.Cells.Item(42, 11).Value = "Sin especificar"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 180
$ws.Cells.Item(42, 14).Value = 7000
$ws.Cells.Item(42, 15).Value = 7500
$ws.Cells.Item(42, 16).Value = 7250
$ws.Cells.Item(42, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(42, 18).Value = "Perú"
$ws.Cells.Item(42, 19).Value = 1812
$ws.Cells.Item(42, 20).Value = 4
